# Update of league bases (data refresh) - swap values between re-ordered
# match rows. The match "id" / rank in column A stays tied to its row
# position, but all the match data (columns B:AD) moves with the match
# record itself, so pairs of rows effectively exchange their B:AD content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($RowA, $RowB)

    $addrA = "B{0}:AD{0}" -f $RowA
    $addrB = "B{0}:AD{0}" -f $RowB

    $rangeA = $ws.Range($addrA)
    $rangeB = $ws.Range($addrB)

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-RowData 134 137
Swap-RowData 135 136
Swap-RowData 139 140
Swap-RowData 142 145
Swap-RowData 143 144
Swap-RowData 255 256
